# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (see header in G1). Update rows 2-13 with
# the newly regenerated K values (previously Strike# derived values).
$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 3
    11 = 1
    12 = 1
    13 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
